# Automatische test-sync: 2025-06-19 19:44:30
# Appends the new "Sollicitatie / Vacature" mail-log entry to the Logs sheet,
# bumps the matching Dashboard category count, and extends the dependent
# conditional-formatting ranges and chart series references to cover the
# new row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Logs sheet: append row 41 with the new entry.
# ---------------------------------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A41").Value = "Sollicitatie marketingfunctie"
$logs.Range("B41").Value = "mailmind.test@zohomail.eu"
$logs.Range("C41").Value = "Hierbij solliciteer ik voor de functie van marketeer. Zie bijlage voor CV."
$logs.Range("D41").Value = "Sollicitatie / Vacature"
$logs.Range("F41").Value = "2025-06-19 19:44:29"
$logs.Range("G41").Value = "Nee"

# Extend the conditional formatting ranges (D2:D40 -> D2:D41, G2:G40 -> G2:G41)
# so the newly added row picks up the same category/answer colour rules.
$dFormatConditions = $logs.Range("D2:D40").FormatConditions
for ($i = 1; $i -le $dFormatConditions.Count; $i++) {
    $dFormatConditions.Item($i).ModifyAppliesToRange($logs.Range("D2:D41"))
}

$gFormatConditions = $logs.Range("G2:G40").FormatConditions
for ($i = 1; $i -le $gFormatConditions.Count; $i++) {
    $gFormatConditions.Item($i).ModifyAppliesToRange($logs.Range("G2:G41"))
}

# ---------------------------------------------------------------------
# 2. Dashboard sheet: add the new category count row 11.
# ---------------------------------------------------------------------
$dashboard = $wb.Worksheets.Item("Dashboard")

$dashboard.Range("A11").Value = "Sollicitatie / Vacature"
$dashboard.Range("B11").Value = 1

# ---------------------------------------------------------------------
# 3. Chart on the Dashboard sheet: extend series cat/val refs to row 11.
# ---------------------------------------------------------------------
$chartObj = $dashboard.ChartObjects(1)
$series = $chartObj.Chart.SeriesCollection(1)
$series.Formula = "=SERIES(Dashboard!`$B`$1,Dashboard!`$A`$2:`$A`$11,Dashboard!`$B`$2:`$B`$11,1)"
